# Fruta / hortaliza, semanal
#
# Insert a new weekly observation row (row 20), pushing the existing
# rows 20-42 down to 21-43, and bump the date/volume of the top row
# (row 19) to reflect the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shifts old rows 20..42 down to 21..43
# and copies row 19's formatting (style) into the freshly inserted row.
$ws.Rows.Item(20).Insert()

# The new row 20 carries forward the data that used to live in row 19.
$ws.Range("A19:R19").Copy()
$ws.Range("A20").PasteSpecial()

# Row 19 itself becomes the newest observation: later date, higher volume.
$ws.Range("D19").Value = 44797
$ws.Range("J19").Value = 60
